$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "queryEntityMapToClickhouseTable"

# --- Bring over header/data-row formatting (fill + border) from an existing sheet ---
$srcSheet = $wb.Worksheets.Item("getConceptModelDataByCondition")
$srcSheet.Range("A1:M2").Copy()
$newSheet.Range("A1:M2").PasteSpecial(-4122)  # xlPasteFormats

$srcSheet.Range("A2:B2").Copy()
$newSheet.Range("A3:B6").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Header row ---
$newSheet.Range("A1").Value = "test-id"
$newSheet.Range("B1").Value = "description"
$newSheet.Range("C1").Value = "condition"
$newSheet.Range("D1").Value = "domainName"
$newSheet.Range("E1").Value = "fields"
$newSheet.Range("F1").Value = "name"
$newSheet.Range("G1").Value = "order"
$newSheet.Range("H1").Value = "pageIndex"
$newSheet.Range("I1").Value = "pageSize"
$newSheet.Range("J1").Value = "timeout"
$newSheet.Range("K1").Value = "rspStatus"
$newSheet.Range("L1").Value = "rspCode"
$newSheet.Range("M1").Value = "rspMessage"

# --- Row 2 ---
$newSheet.Range("A2").Value = "iot-connector-enittymaptoclickhousetable-1"
$newSheet.Range("B2").Value = "good request, data retrieved"
$newSheet.Range("F2").Value = "ClickhouseDriverSensor"

# --- Row 3 ---
$newSheet.Range("A3").Value = "iot-connector-enittymaptoclickhousetable-2"
$newSheet.Range("B3").Value = "good request, data retrieved"
$newSheet.Range("F3").Value = "ClickhouseDriverSensor"
$newSheet.Range("H3").Value = 1
$newSheet.Range("I3").Value = 2

# --- Row 4 ---
$newSheet.Range("A4").Value = "iot-connector-enittymaptoclickhousetable-3"
$newSheet.Range("B4").Value = "good request, data retrieved"
$newSheet.Range("F4").Value = "ClickhouseDriverSensor"
$newSheet.Range("G4").Value = "value"

# --- Row 5 ---
$newSheet.Range("A5").Value = "iot-connector-enittymaptoclickhousetable-4"
$newSheet.Range("B5").Value = "good request, data retrieved"
$newSheet.Range("C5").Value = "name='sensorA'"
$newSheet.Range("F5").Value = "ClickhouseDriverSensor"

# --- Row 6 ---
$newSheet.Range("A6").Value = "iot-connector-enittymaptoclickhousetable-5"
$newSheet.Range("B6").Value = "good request, data retrieved"
$newSheet.Range("E6").Value = "timestamp"
$newSheet.Range("F6").Value = "ClickhouseDriverSensor"

# --- Column widths (match target) ---
$newSheet.Columns.Item(1).ColumnWidth = 62.5546875
$newSheet.Columns.Item(2).ColumnWidth = 27.6640625
$newSheet.Columns.Item(3).ColumnWidth = 26
$newSheet.Columns.Item(4).ColumnWidth = 16.5546875
$newSheet.Columns.Item(5).ColumnWidth = 23.109375
$newSheet.Columns.Item(6).ColumnWidth = 28.6640625

# --- View: selected cell + scroll position ---
$newSheet.Range("C5").Select()

$wb.Windows.Item(1).ScrollColumn = 2
